$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "role" data currently in row 7 (columns B-J) down to a new row 11,
# and clear those cells from row 7.

$ws.Range("B11").Value = $ws.Range("B7").Value2
$ws.Range("C11").Value = $ws.Range("C7").Value2
$ws.Range("D11").Value = $ws.Range("D7").Value2
$ws.Range("E11").Value = $ws.Range("E7").Value2
$ws.Range("F11").Value = $ws.Range("F7").Value2
$ws.Range("G11").Value = $ws.Range("G7").Value2
$ws.Range("H11").Value = $ws.Range("H7").Value2
$ws.Range("I11").Value = $ws.Range("I7").Value2
$ws.Range("J11").Value = $ws.Range("J7").Value2

$ws.Range("B7:J7").ClearContents()
